$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update rq_excl_id description (row 6 / column D) to mention duplicate ids
#    and grow the row height to fit the extra lines of text.
# ---------------------------------------------------------------------------
$ws.Rows.Item(6).RowHeight = 57
$ws.Cells.Item(6,4).Value = "VRM2 shall enable exclusion of requirements from diagram based on <id> and in case of duplicate <id>s the key created from id:version.`nComment: exclusion is typically done interactively, where the relevant identifier is picked up programatically with a right-click on relevant requirement.`n"

# ---------------------------------------------------------------------------
# 2. Finish off row 70 (previously only had a stray D70 value), then add the
#    two brand new requirement rows (71, 72) covering duplicate handling.
# ---------------------------------------------------------------------------
$ws.Rows.Item(70).RowHeight = 24
$ws.Cells.Item(70,1).Value = "rq_dup_req"
$ws.Cells.Item(70,2).Value = 1
$ws.Cells.Item(70,3).Value = "approved"
$ws.Cells.Item(70,4).Value = "VRM2 shall handle duplicate ids and calculate links based on id+version`nIf no matching version is found for a fulfilledby or coverage relation, it is unspecified which duplicate is linked to."
$ws.Cells.Item(70,5).Value = "impl;test"

$ws.Rows.Item(71).RowHeight = 13
$ws.Cells.Item(71,1).Value = "rq_dup_req_display"
$ws.Cells.Item(71,2).Value = 1
$ws.Cells.Item(71,3).Value = "approved"
$ws.Cells.Item(71,4).Value = "VRM2 shall handle group duplicate requirements together visually"
$ws.Cells.Item(71,5).Value = "impl;test"

$ws.Rows.Item(72).RowHeight = 13
$ws.Cells.Item(72,1).Value = "rq_dup_same_version"
$ws.Cells.Item(72,2).Value = 1
$ws.Cells.Item(72,3).Value = "approved"
$ws.Cells.Item(72,4).Value = "VRM2 shall log an issue when a non-unique id+version pair is detected"
$ws.Cells.Item(72,5).Value = "impl;test"

# ---------------------------------------------------------------------------
# 3. Re-apply the autofilter over the grown A1:G72 range, filtering column D
#    (index 4) down to the six requirements that stay visible. Clearing the
#    old filter first forces the filter range itself to grow to G72 too.
# ---------------------------------------------------------------------------
$ws.AutoFilterMode = $false
$filterValues = @(
    "VRM2 shall calculate the set of shown requirements based on reachability from selected requirements through coverage relations.`nExcluded ids shall stop the graph traversal and thus limit the set of shown requirements.`nExcluded doctypes shall stop the graph traversal and thus limit the set of shown requirements.",
    "VRM2 shall display a legend containing input file name(s), selection criteria, excluded ids and `u2018safety`u2019 rules, i.e. the information needed to reproduce the diagram.",
    "VRM2 shall enable exclusion of requirements from diagram based on <id>",
    "VRM2 shall enable exclusion of requirements from diagram based on doctypes",
    "VRM2 shall provide a mechanism to exclude `u2018rejected`u2019 requirements from the diagram",
    "VRM2 shall provide an option to exclude a specobject from the context menu"
)
$ws.Range("A1:G72").AutoFilter(4, $filterValues, 7)

# ---------------------------------------------------------------------------
# 4. Fix up the defined name driving _FilterDatabase so it matches the new
#    autofilter range (the engine does not refresh it automatically).
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name() -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$G`$72"
    }
}

# ---------------------------------------------------------------------------
# 5. Move the active selection, matching the saved view state in the diff.
# ---------------------------------------------------------------------------
$ws.Range("D60").Select()
